$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 385, shifting existing rows 385..449 down to 386..450
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new record's data
$ws.Cells.Item(385, 1).Value = 4
$ws.Cells.Item(385, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(385, 3).Value = 'Los Lagos'
$ws.Cells.Item(385, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(385, 5).Value = 10
$ws.Cells.Item(385, 6).Value = 'Fruta'
$ws.Cells.Item(385, 7).Value = 100108
$ws.Cells.Item(385, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(385, 9).Value = 100108002
$ws.Cells.Item(385, 10).Value = 'Mango'
$ws.Cells.Item(385, 11).Value = 'Sin especificar'
$ws.Cells.Item(385, 12).Value = 'Primera'
$ws.Cells.Item(385, 13).Value = 100
$ws.Cells.Item(385, 14).Value = 13500
$ws.Cells.Item(385, 15).Value = 13500
$ws.Cells.Item(385, 16).Value = 13500
$ws.Cells.Item(385, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(385, 18).Value = 'Brasil'
$ws.Cells.Item(385, 19).Value = 3375
$ws.Cells.Item(385, 20).Value = 4
